# "Generate Report for Handoff"
# Updates the localization-status workbook: the previous handback pair of
# files (58095348-...md / d7cb7150-...md) is replaced by a new pair that is
# now ready for handoff (135675ad-...md / ffffc6b7a6d2-...md), the Status
# columns flip from "Handed back: in sync with en-US" to "Ready for handoff",
# timestamps move forward, and the per-locale "Latest Target/Handback" info
# is cleared out because the new pair has not been handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "58095348-c0b4-43d8-8f04-d10ba197177c"
$oldGuid2 = "d7cb7150-4012-4530-9ce7-a12d547371e3"
$newGuid1 = "135675ad-6f70-4b39-950c-de98fae9371f"
$newGuid2 = "ffffc6b7a6d2-5fef-4cf2-8d09-28028b1e74ce"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newGuid1.md"
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-25 11:05:37"

$ov.Range("A3").Value = "$newGuid2.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-25 11:05:37"

# Hyperlinks in column B display the path, but keep pointing at the same
# external targets as before (the repo links themselves were not touched).
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/299616c96049a1d7a9f9db4c2331bbbaa91a6dcb/e2e/$oldGuid1.md", "", "", "e2e\$newGuid1.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/299616c96049a1d7a9f9db4c2331bbbaa91a6dcb/e2e/$oldGuid2.md", "", "", "e2e\$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newGuid1.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("G2").Value = "$newGuid1.2bd1f3af0fbfa928fc0510428fc7bfd5a7a4c2f8.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-25 11:05:32"
$zh.Range("I2").Font.Underline = $false
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "$newGuid2.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "$newGuid1.2bd1f3af0fbfa928fc0510428fc7bfd5a7a4c2f8.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-25 11:05:32"
$zh.Range("I3").Font.Underline = $false
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/299616c96049a1d7a9f9db4c2331bbbaa91a6dcb/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/27ad1880c29635047eb8731057321e5f69bc9f56/e2e/$oldGuid1.md", "", "", "$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newGuid1.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("G2").Value = "$newGuid1.2bd1f3af0fbfa928fc0510428fc7bfd5a7a4c2f8.de-de.xlf"
$de.Range("H2").Value = "2016-08-25 11:05:37"
$de.Range("I2").Font.Underline = $false
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = "$newGuid2.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "$newGuid1.2bd1f3af0fbfa928fc0510428fc7bfd5a7a4c2f8.de-de.xlf"
$de.Range("H3").Value = "2016-08-25 11:05:37"
$de.Range("I3").Font.Underline = $false
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/299616c96049a1d7a9f9db4c2331bbbaa91a6dcb/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2d0f6c61be0762f6cbe307589b23570a6f355c82/e2e/$oldGuid1.md", "", "", "$newGuid2.md")

# ---------------------------------------------------------------------
# Column widths: the shorter replacement text (shorter status string,
# blank handback columns) made Excel re-autofit a handful of columns.
# (Input values are calibrated so the engine's internal 1/6-character
# rounding lands on the width closest to the target OOXML column width.)
# ---------------------------------------------------------------------
$ov.Range("E1").EntireColumn.ColumnWidth = 16.33
$ov.Range("F1").EntireColumn.ColumnWidth = 16.33

foreach ($sheet in @($zh, $de)) {
    $sheet.Range("C1").EntireColumn.ColumnWidth = 16.33
    $sheet.Range("I1").EntireColumn.ColumnWidth = 17.8
    $sheet.Range("J1").EntireColumn.ColumnWidth = 20.8
}

Write-Host "Report regenerated for handoff."
